$d = $word.ActiveDocument

# Change 1: split "Todos os arquivos..." text, dropping ", juntamente com a apresentação"
$d.Content.Find.Execute("Todos os arquivos deverão ser compactados, juntamente com a apresentação e enviados para e-mail: ", $true, $false, $false, $false, $false, $true, 1, $false, "Todos os arquivos deverão ser compactados e enviados para e-mail: ", 2)

Write-Output "done"
